# "Added 1 more charge in Shipment Creation"
# Updates the ShipmentTracking numbers (column P) for the two new/changed
# charge rows, replacing the old tracking numbers with the new ones.
#
# Assigning a purely-numeric string straight to Range.Value turns the cell
# into a Number cell. The source values in this sheet are stored as TEXT
# (shared strings), so we round-trip through a TEXT() formula and then
# Copy / PasteSpecial (values only) to collapse it back down to a literal
# string cell without dragging a new NumberFormat/style onto the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Formula = '=TEXT(' + $text + ',"0")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue $ws.Cells.Item(3, 16) "320018720183"
Set-TextValue $ws.Cells.Item(5, 16) "320018720210"
